$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 163 ("Bruselas (repollito)"
# at Vega Modelo de Temuco, 29-May-2023), pushing the previous rows 163-169
# down to 164-170. Insert a fresh row at 163 to shift everything down.
$ws.Rows(163).Insert()

# Populate the newly inserted row 163 with the new record's data.
$ws.Cells.Item(163, 1).Value = 10
$ws.Cells.Item(163, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(163, 3).Value = "La Araucanía"
$ws.Cells.Item(163, 4).Value = 45075
$ws.Cells.Item(163, 5).Value = 9
$ws.Cells.Item(163, 6).Value = 100112035
$ws.Cells.Item(163, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 35
$ws.Cells.Item(163, 11).Value = 28000
$ws.Cells.Item(163, 12).Value = 28000
$ws.Cells.Item(163, 13).Value = 28000
$ws.Cells.Item(163, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(163, 15).Value = "Región Metropolitana"
$ws.Cells.Item(163, 16).Value = 1867
$ws.Cells.Item(163, 17).Value = 15
$ws.Cells.Item(163, 18).Value = "Hortaliza"

# Match the date cell format used by the rest of column D (row above),
# so the inserted row carries style index 2, same as its siblings.
$ws.Cells.Item(163, 4).NumberFormat = $ws.Cells.Item(164, 4).NumberFormat()
